$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new "source_name" column before the existing column B (junction_type),
# shifting the original B:H columns to C:I.
$ws.Columns("B").Insert()

$ws.Range("B1").Value = "source_name"

# Fill in the new data row (row 2). Enter the carried-over junction data
# (columns D:I) first, then the new source columns (B2/C2) last, matching
# the original authoring order.
$ws.Range("A2").Value = 1
$ws.Range("D2").Value = "GGCGCU"
$ws.Range("E2").Value = "G"
$ws.Range("F2").Value = "GCCUGCAGGC"
$ws.Range("G2").Value = "UA"
$ws.Range("H2").Value = "GACGUG"
$ws.Range("I2").Value = "A"

$ws.Range("C2").Value = "Family_A"
$ws.Range("B2").Value = "16S H20-H21-H22"

# Best-fit column widths (as Excel would after auto-fitting the new data).
$ws.Columns("A").ColumnWidth = 8
$ws.Columns("B").ColumnWidth = 44/3
$ws.Columns("C").ColumnWidth = 34/3
$ws.Columns("D").ColumnWidth = 23/3
$ws.Columns("E").ColumnWidth = 53/6
$ws.Columns("F").ColumnWidth = 73/6
$ws.Columns("G").ColumnWidth = 53/6
$ws.Columns("H").ColumnWidth = 47/6
$ws.Columns("I").ColumnWidth = 26/3

$null = $ws.Range("L8").Select()
